$d = $word.ActiveDocument

# --- Step 1: turn the " m: 2.myTemplate() " field into literal text "{m: 2.myTemplate()}" ---
$f = $d.Fields.Item(1)
$codeRange = $f.Code
$beginPos = $codeRange.Start - 1

$newText = "{m: 2.myTemplate()}"
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($beginPos, $beginPos)
$insertionPoint.InsertXML($xmlFrag)

# the fldChar/instrText/fldChar run-triplet is still right after the text we just inserted;
# remove it (Field.Delete removes all three runs that make up the field).
$f2 = $d.Fields.Item(1)
$f2.Delete()

# --- Step 2: prefix the error message with "    <---" ---
$prefixLen = $newText.Length
$errStart = $beginPos + $prefixLen
$errCharRange = $d.Range($errStart, $errStart + 1)
$errCharRange.Text = "    <---" + $errCharRange.Text
